$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Both edits insert manual line breaks (<w:br/>) inside a single run's text,
# splitting one long <w:t> into several <w:t>/<w:br/> runs. We locate each
# target paragraph by its distinctive leading text, then rewrite its content
# with Range.InsertXML so we get exact control over run/break placement and
# xml:space="preserve" (needed on runs that end in a trailing space).
# ---------------------------------------------------------------------------

$pPrograma = $null
$pBibliografia = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("Formação e trabalho em equipes e Comunicação")) {
        $pPrograma = $p
    }
    elseif ($t.StartsWith("Gestão de Negócios: Visões e dimensões empresariais")) {
        $pBibliografia = $p
    }
}

if ($pPrograma -eq $null) {
    throw "Could not locate the 'Programa' paragraph to edit"
}
if ($pBibliografia -eq $null) {
    throw "Could not locate the 'Bibliografia' paragraph to edit"
}

# --- Edit 1: "Programa" section ---------------------------------------------
# Break the paragraph right after "(mínimo 3);" and before "Inovação
# Sistemática", turning it into two runs joined by a <w:br/>.
$programaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:r>' +
    '<w:t>Formação e trabalho em equipes e Comunicação – o desenvolvimento das habilidades essenciais para o trabalho em equipes ocorrerá por meio do trabalho em equipes e de reuniões e visitas didáticas realiadas na empresa (mínimo 3);</w:t>' +
    '<w:br/>' +
    '<w:t>Inovação Sistemática – desenvolvimento de soluções inovadoras, sistematização e características; Legislação - noções da legislação aplicada à ação empresarial; Gerenciamento de Projetos e Cronograma – Metodologias e esquematizações necessárias com os elementos gerenciais; Identificação de Problemas – sistematização de ações para a localização de causas; Formulação do Projeto – apresentação dos aspectos gerenciais necessários ao desenvolvimento do projeto, Plano de gestão, Estrutura Analítica do Projeto (EAP) etc; Especificação de Problemas – sistematização dos problemas dentro das áreas de conhecimento; Análise do Conhecimento disponível, Avaliação e Tomada de Decisão; Elaboração de relatórios – formatação dentro das normas ABNT; Apresentação de Projetos.</w:t>' +
    '</w:r>' +
    '</w:p>'
$pPrograma.Range.InsertXML($programaXml)

# --- Edit 2: "Bibliografia" section -----------------------------------------
# Break the single run into many runs joined by <w:br/> (doubled between
# distinct references), one per citation/field line.
$bibliografiaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:r>' +
    '<w:t xml:space="preserve">Gestão de Negócios: Visões e dimensões empresariais da o Organização. </w:t>' +
    '<w:br/>' +
    '<w:t>Autores: Cruz Jr, J.B., Rocha, J.A.O. e Tachizawa, T.</w:t>' +
    '<w:br/>' +
    '<w:t>Editora: ATLAS</w:t>' +
    '<w:br/>' +
    '<w:br/>' +
    '<w:t>Gestão Empresarial - de Taylor aos nossos dias</w:t>' +
    '<w:br/>' +
    '<w:t xml:space="preserve">Autores: Pereira, M. I. , Autor: Ferreira, A. A. e Reis, A.C. F </w:t>' +
    '<w:br/>' +
    '<w:t>Editora: THOMSON PIONEIRA</w:t>' +
    '<w:br/>' +
    '<w:br/>' +
    '<w:t>Baron e Shane: Empreendedorismo: uma visão do processo (EVP), Ed. Thomson, 2006</w:t>' +
    '<w:br/>' +
    '<w:br/>' +
    '<w:t>Textos fornecidos pelo professor da disciplina</w:t>' +
    '<w:br/>' +
    '<w:t>Artigos extraídos de revistas especializadas na área de gestão e produção.</w:t>' +
    '</w:r>' +
    '</w:p>'
$pBibliografia.Range.InsertXML($bibliografiaXml)
